# ABML de tipo de Examen - corrección de nombres/apellidos de usuarios
# (Profesores y Alumnos de prueba) en la hoja "2.Usuarios", y ajuste de
# la vista/selección activa tras la edición.

$wb = $excel.ActiveWorkbook

$wsUsuarios = $wb.Worksheets.Item("2.Usuarios")
$wsMateria  = $wb.Worksheets.Item("3.Materia-Curso")

# --- Profesores (filas 7 y 8) ---
$wsUsuarios.Range("B7").Value = "Ayelén"
$wsUsuarios.Range("C7").Value = "Flores"

$wsUsuarios.Range("B8").Value = "Silvia"
$wsUsuarios.Range("C8").Value = "Malloti"

# --- Alumnos (filas 9 a 12) ---
$wsUsuarios.Range("B9").Value = "Clara"
$wsUsuarios.Range("C9").Value = "Del Valle"

$wsUsuarios.Range("B10").Value = "Jaime"
$wsUsuarios.Range("C10").Value = "Trueba"

$wsUsuarios.Range("B11").Value = "Alba"
$wsUsuarios.Range("C11").Value = "De Satigny"

$wsUsuarios.Range("B12").Value = "Pedro"
$wsUsuarios.Range("C12").Value = "Tercero García"

# Deja la hoja de Materia-Curso con el rango revisado seleccionado.
$wsMateria.Activate()
$wsMateria.Range("F2:F5").Select()

# La hoja de Usuarios queda como hoja activa, con el cursor en B7.
$wsUsuarios.Activate()
$wsUsuarios.Range("B7").Select()
